$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - Mobile Numeric Keypad Problem
$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "Mobile Numeric Keypad Problem"
$ws.Range("C6").Value = "DP"
$ws.Range("D6").Value = "Array"
$ws.Range("E6").Value = "easy"
$ws.Range("F6").Value = "GeeksForGeeks"

# Row 7 - Total number of non-decreasing numbers with n digits
$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Total number of non-decreasing numbers with n digits"
$ws.Range("C7").Value = "DP"
$ws.Range("D7").Value = "Array"
$ws.Range("E7").Value = "easy"
$ws.Range("F7").Value = "GeeksForGeeks"

# Update the active selection to D14, matching the saved view state
$ws.Range("D14").Select()
